$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7: policy end date shifted; the travel_limit end-marker (W7) flips 1 -> 0
# and the day-count denominator (AA7) drops from 13 to 12. ---
$ws.Range("W7").Value = 0
$ws.Range("AA7").Value = 12

# --- Recompute the day-weighted AA column (rows 25-221) against the new 12-day span. ---
$ws.Range("AA25:AA31").Value = 0.08333333333333333
$ws.Range("AA32:AA39").Value = 0.1666666666666667
$ws.Range("AA40:AA63").Value = 0.7499999999916667
$ws.Range("AA64:AA70").Value = 0.5166666666583333
$ws.Range("AA71:AA84").Value = 0.3194444444416667
$ws.Range("AA85:AA93").Value = 0.25
$ws.Range("AA94:AA221").Value = 0.1666666666666667

# --- Append 12 new daily rows (222-233) for 9/30/2020 .. 10/11/2020, mirroring the
#     "no new restriction that day" pattern already used by the preceding rows. ---
$newDates = @("9/30/2020", "10/1/2020", "10/2/2020", "10/3/2020", "10/4/2020", "10/5/2020", "10/6/2020", "10/7/2020", "10/8/2020", "10/9/2020", "10/10/2020", "10/11/2020")
$rowVals = @(0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0)
$r = 222
foreach ($d in $newDates) {
    $ws.Range("A" + $r).NumberFormat = "@"
    $ws.Range("A" + $r).Value = $d
    $ws.Range("A221").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $col = 2
    foreach ($v in $rowVals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
    $ws.Cells.Item($r, 27).Value = 0.1666666666666667

    $r++
}
